$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.089.10'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.654.48'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '''217.69'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '''0.5250'
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '''0.2608'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("D9").Value = '''0.06348'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").Value = '''20.35'
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("D11").Value = '''0.07805'
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").Value = '1.659.89'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '0.0₅8209'
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").Value = '''65.35'
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").Value = '26.105.61'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '''1.002'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '''4.582'
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").Value = '''191.19'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '''6.016'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = '''142.08'
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("D25").Value = '''0.1238'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Value = '''7.255'
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").Value = '''16.10'
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").Value = '''1.428'
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("D29").Value = '''0.05900'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").Value = '''3.510'
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = '''3.251'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''1.588'
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").Value = '''0.9501'
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("D35").Value = '''2.782'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("D37").Value = '''0.5688'
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").Value = '''0.01619'
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("D39").Value = '''5.809'
$ws.Range("E39").Value = '  -2.73%  '
$ws.Range("D40").Value = '''0.8494'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D42").Value = '''102.84'
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("D43").Value = '1.026.84'
$ws.Range("D44").Value = '1.797.18'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '''0.4304'
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").Value = '''0.05163'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").Value = '''7.833'
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("D51").Value = '''0.09701'
$ws.Range("E51").Value = '  -0.18%  '
